# Fill in the "Pass/Fail" (column F) results and the missing "Actual Result"
# note for the login-while-logged-out test case, per the updated test-case
# spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1 (row 2): while logged out -> returned to login page => Fail
$ws.Range("F2").Value = "Fail"

# Step 2 (row 3): login as appropriate role -> redirected to dashboard => Pass
$ws.Range("F3").Value = "Pass"

# Step 3 (row 4): go to "Create Employee" -> redirected to Create Employee page => Pass
$ws.Range("F4").Value = "Pass"

# Step 4 (row 5): submit employee info -> should add employee to database,
# but nothing happened (no error message) => Fail
$ws.Range("E5").Value = "Nothing happened. No error message."
$ws.Range("F5").Value = "Fail"

# Update the view so column C is the left-most visible column and E4 is
# the active selection.
$ws.Range("E4").Select()
